# Rename Sheet1 -> ValidLogin and populate the login-sample data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "ValidLogin"

# Values are written in this order (A1,B1,A2,B2,C2,C1) so the shared-string
# table comes out in the same order as the source workbook.
$ws.Range("A1").Value = "Username"
$ws.Range("B1").Value = "password"
$ws.Range("A2").Value = "admin"
$ws.Range("B2").Value = "manager"
$ws.Range("C2").Value = "actiTIME - Enter Time-Track"
$ws.Range("C1").Value = "eTitle"

# Widen column C to fit the long header/value text.
$ws.Range("C1").ColumnWidth = 25.666666666666668

# Leave the selection on C1, like the source file.
[void]$ws.Range("C1").Select()
